$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data as scraped on Sat Aug 31 06:51:29 UTC 2024

$ws.Range("D2").Value = "59.201.75"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.524.83"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'537.16"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "'137.96"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "2.523.74"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "'0.350"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "2.974.41"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "'23.16"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "59.144.63"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "2.527.64"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'11.13"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "'325.91"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'5.99"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").Value = "'66.18"
$ws.Range("E24").Value = "  +6.29%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").Value = "'6.74"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0776"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = "  +5.66%  "
$ws.Range("D33").Value = "'162.94"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.47"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").Value = "'0.821"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'287.38"
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "'132.42"
$ws.Range("E44").Value = "  +8.24%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "'0.611"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'0.0932"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").Value = "'17.43"
$ws.Range("E51").Value = "  -1.47%  "
